$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''330.42'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''3.26%'
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''41.16'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''4.67%'
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = '''5.694'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '''-3.24%'
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''0.08172'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''2.03%'
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''2.081'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''9.88%'
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = '''8.731'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''1.13%'
$ws.Range("E7").Style = "Normal"
$ws.Range("B8").Value = 'GateToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D8").Value = '''4.539'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''-0.73%'
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = 'BTSEToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D9").Value = '''2.962'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''0.38%'
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = 'MXToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D10").Value = '''0.9257'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''-1.02%'
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").Value = '''0.1263'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''1.02%'
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = 'WazirX'
$ws.Range("C12").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D12").Value = '''0.1963'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''0.73%'
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = 'MandalaExchangeToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D13").Value = '''0.09451'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''3.66%'
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").Value = '''0.03688'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''5.09%'
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").Value = '''0.1055'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''10.21%'
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").Value = '''0.001297'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''1.18%'
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D17").Value = '''0.04427'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''-0.85%'
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = 'TigerCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D18").Value = '''0.006149'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''0.75%'
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = 'LEO'
$ws.Range("C19").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D19").Value = '''3.414'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''1.76%'
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = 'BitpandaEcosystemToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D20").Value = '''0.3485'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''-1.46%'
$ws.Range("E20").Style = "Normal"
$ws.Range("B21").Value = 'MCDex'
$ws.Range("C21").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D21").Value = '''8.308'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''-5.06%'
$ws.Range("E21").Style = "Normal"
$ws.Range("B22").Value = 'ProBitToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D22").Value = '''0.1380'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''-3.47%'
$ws.Range("E22").Style = "Normal"
$ws.Range("B23").Value = 'ZBToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D23").Value = '''0.2652'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''9.95%'
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''0.001269'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''0.38%'
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''0.004345'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''-1.63%'
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''0.0001182'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''3.45%'
$ws.Range("E26").Style = "Normal"
$ws.Range("D39").Value = '''0.02812'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''17.29%'
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''0.05481'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''5.95%'
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''0.007679'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''3.34%'
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''0.009432'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''3.00%'
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = '''1.19%'
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = '''0.42%'
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''0.01178'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''5.49%'
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''0.00006869'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''1.86%'
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = '''-0.05%'
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = '''60.43%'
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = '''0.003161'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''5.00%'
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''0.00002103'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''-0.05%'
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''0.0002003'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''-0.05%'
$ws.Range("E51").Style = "Normal"
